$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Replace-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.InsertXML($innerXml)
}

# ---------------------------------------------------------------------------
# NOTE: edits are applied from the bottom of the document upward so that
# paragraph indices captured above an edit point never shift underneath us.
# ---------------------------------------------------------------------------

# --- Experience > Freelance Projects (paragraph 39): new bullet list -------
$freelance = "<w:p $wNs><w:pPr><w:pStyle w:val='BodyText'/></w:pPr>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'>Freelance Projects</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>- Built responsive websites using WordPress and Elementor</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>- Worked autonomously with clients to deliver solutions on time</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>- Managed requirements, execution, and final delivery independently</w:t></w:r>" +
    "</w:p>"
Replace-ParagraphXml 39 $freelance

# --- Experience > second entry (paragraph 38): now ESSABHY (BodyText) ------
$essabhy = "<w:p $wNs><w:pPr><w:pStyle w:val='BodyText'/></w:pPr>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'>ESSABHY</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>– Tech Company</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>- Developed WordPress-based websites and business tools</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>- Participated in team-based project delivery and client support</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>- Gained exposure to production environments and client-facing solutions</w:t></w:r>" +
    "</w:p>"
Replace-ParagraphXml 38 $essabhy

# --- Experience > first entry (paragraph 37): now ITNERA (FirstParagraph) --
$itnera = "<w:p $wNs><w:pPr><w:pStyle w:val='FirstParagraph'/></w:pPr>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'>ITNERA</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>– Tech Company</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>- Led and supervised teams delivering websites, desktop applications, and IT solutions</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>- Coordinated project execution from requirements to deployment</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>- Installed, configured, and deployed digital solutions and technical equipment for clients</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>- Contributed to e-commerce platforms, SaaS-style tools, and business applications</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>- Provided technical support and ensured system reliability for client environments</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>- Applied backend and system knowledge to improve workflows and solution quality</w:t></w:r>" +
    "</w:p>"
Replace-ParagraphXml 37 $itnera

# --- Skills (paragraph 6): split into Languages/WebBasics/C++, Backend&Tools, Domains&Focus
$skills1 = "<w:p $wNs><w:pPr><w:pStyle w:val='FirstParagraph'/></w:pPr>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'>Languages:</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>C, C++, JavaScript</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'>Web Basics:</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>HTML, CSS, PHP (foundational knowledge)</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'>C++ Libraries:</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>STL, Networking and system libraries</w:t></w:r>" +
    "</w:p>"
$skills2 = "<w:p $wNs><w:pPr><w:pStyle w:val='BodyText'/></w:pPr>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'>Backend &amp; Tools:</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>Node.js, NestJS, Docker, Docker Compose, PostgreSQL, Linux, Bash, Git, Vue.js</w:t></w:r>" +
    "</w:p>"
$skills3 = "<w:p $wNs><w:pPr><w:pStyle w:val='BodyText'/></w:pPr>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'>Domains &amp; Focus:</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>Performance programming, Security, DevOps, e-commerce platforms, SaaS systems,</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>business-oriented applications</w:t></w:r>" +
    "</w:p>"
$skillsAll = $skills1 + $skills2 + $skills3
Replace-ParagraphXml 6 $skillsAll

# --- Contact block (paragraph 3): new summary paragraph + reordered contact info
$summary = "<w:p $wNs><w:pPr><w:pStyle w:val='BodyText'/></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>Backend / Full-Stack Developer (Intern / Junior) with strong foundations in</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>systems programming, backend development, and security. Experienced in team</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>collaboration and delivering real-world solutions, including e-commerce and</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>business platforms. Autonomous, reliable with deadlines, and quick to learn in</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>hands-on environments.</w:t></w:r>" +
    "</w:p>"
$contact = "<w:p $wNs><w:pPr><w:pStyle w:val='BodyText'/></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>Location: Morocco</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>Email: medattiq@gmail.com</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>GitHub: https://github.com/simbaattiq</w:t></w:r>" +
    "</w:p>"
$introAll = $summary + $contact
Replace-ParagraphXml 3 $introAll

Write-Output "done"
